$d = $word.ActiveDocument

$replacements = @(
    @("2025-11-30 Sunday", "2025-12-01 Monday"),
    @("814÷6=", "155÷3="),
    @("411÷9=", "345÷5="),
    @("470÷6=", "473÷3="),
    @("542÷6=", "365÷6="),
    @("302÷5=", "972÷8="),
    @("576÷6=", "253÷3="),
    @("130÷3=", "489÷4="),
    @("459÷7=", "653÷8="),
    @("125÷4=", "896÷7="),
    @("782÷8=", "186÷5="),
    @("493÷9=", "692÷4="),
    @("895÷5=", "416÷5="),
    @("162÷2=", "420÷8="),
    @("703÷6=", "278÷9="),
    @("833÷7=", "884÷5="),
    @("541÷8=", "978÷6="),
    @("627÷3=", "178÷9="),
    @("992÷2=", "750÷7="),
    @("845÷3=", "351÷2="),
    @("562÷6=", "868÷6="),
    @("859÷5=", "800÷6="),
    @("710÷4=", "354÷4="),
    @("320÷6=", "581÷9="),
    @("644÷2=", "269÷9="),
    @("940÷2=", "462÷5=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
